$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8066960573196411
$ws.Range("B1").Value = 0.6395947933197021
$ws.Range("C1").Value = 2.206688404083252
$ws.Range("D1").Value = 3.396670341491699
$ws.Range("E1").Value = 1.429866313934326
